# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the "Palta" data block
# (rows 958-959), pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 958; everything from the old row 958
# onward shifts down to row 960 onward.
$ws.Rows("958:959").Insert()

# --- New row 958: Hass / Primera -----------------------------------------
$ws.Range("A958").Value = 11
$ws.Range("B958").Value = 'Vega Monumental Concepción'
$ws.Range("C958").Value = 'Bíobío'
$ws.Range("D958").Value = 45132
$ws.Range("E958").Value = 8
$ws.Range("F958").Value = 'Fruta'
$ws.Range("G958").Value = 100106
$ws.Range("H958").Value = 'Oleaginosos'
$ws.Range("I958").Value = 100106002
$ws.Range("J958").Value = 'Palta'
$ws.Range("K958").Value = 'Hass'
$ws.Range("L958").Value = 'Primera'
$ws.Range("M958").Value = 150
$ws.Range("N958").Value = 25000
$ws.Range("O958").Value = 25000
$ws.Range("P958").Value = 25000
$ws.Range("Q958").Value = '$/bandeja 10 kilos'
$ws.Range("R958").Value = 'Perú'
$ws.Range("S958").Value = 2500
$ws.Range("T958").Value = 10

# --- New row 959: Hass / Segunda -----------------------------------------
$ws.Range("A959").Value = 11
$ws.Range("B959").Value = 'Vega Monumental Concepción'
$ws.Range("C959").Value = 'Bíobío'
$ws.Range("D959").Value = 45132
$ws.Range("E959").Value = 8
$ws.Range("F959").Value = 'Fruta'
$ws.Range("G959").Value = 100106
$ws.Range("H959").Value = 'Oleaginosos'
$ws.Range("I959").Value = 100106002
$ws.Range("J959").Value = 'Palta'
$ws.Range("K959").Value = 'Hass'
$ws.Range("L959").Value = 'Segunda'
$ws.Range("M959").Value = 120
$ws.Range("N959").Value = 20000
$ws.Range("O959").Value = 20000
$ws.Range("P959").Value = 20000
$ws.Range("Q959").Value = '$/bandeja 10 kilos'
$ws.Range("R959").Value = 'Perú'
$ws.Range("S959").Value = 2000
$ws.Range("T959").Value = 10
